$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the orphan data row (old row 21), which only had F/G/H totals and
# no idx/name - this shifts all the following rows up by one.
$ws.Rows("21").Delete()

# The old sheet had a two-row header (row 1 + row 2). Collapse it into a
# single header row by deleting the old second header row; this shifts the
# plant data rows up so they start at row 2.
$ws.Rows("2").Delete()

# Rebuild row 1 as one unified header row with new column captions.
$ws.Range("A1:K1").Clear()
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# The numeric/unit columns (F:K) keep the bold-ish "applied font" style that
# the old two-row header used; A:E (idx/idx2/Name/dates) stay plain.
$ws.Range("F1:K1").Font.Size = 9

# Restore the active selection to match the edited document state.
[void]$ws.Range("A20:K20").Select()
